$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set text number format for percentage cells first, so the % strings are
# preserved as literal text instead of being parsed into numeric percentages.
$ws.Range('H2').NumberFormat = "@"
$ws.Range('H5').NumberFormat = "@"
$ws.Range('H7').NumberFormat = "@"
$ws.Range('H8').NumberFormat = "@"
$ws.Range('H9').NumberFormat = "@"
$ws.Range('H11').NumberFormat = "@"
$ws.Range('H12').NumberFormat = "@"
$ws.Range('H13').NumberFormat = "@"
$ws.Range('H15').NumberFormat = "@"
$ws.Range('H16').NumberFormat = "@"
$ws.Range('H17').NumberFormat = "@"
$ws.Range('H21').NumberFormat = "@"
$ws.Range('H24').NumberFormat = "@"
$ws.Range('H25').NumberFormat = "@"
$ws.Range('H26').NumberFormat = "@"
$ws.Range('H29').NumberFormat = "@"
$ws.Range('H31').NumberFormat = "@"
$ws.Range('H32').NumberFormat = "@"
$ws.Range('H35').NumberFormat = "@"
$ws.Range('H36').NumberFormat = "@"
$ws.Range('H42').NumberFormat = "@"
$ws.Range('H46').NumberFormat = "@"

# Apply the updated cell values (latest automatic meteocat.cat extraction).
$ws.Range('E2').Value = '2026-02-09 04:18:22'
$ws.Range('H2').Value = '91%'
$ws.Range('O2').Value = '-3.9 °C'
$ws.Range('E3').Value = '2026-02-09 04:18:25'
$ws.Range('O3').Value = '-6.4 °C'
$ws.Range('E4').Value = '2026-02-09 04:18:27'
$ws.Range('L4').Value = '6.8 km/h - 293º 3:52 TU'
$ws.Range('N4').Value = '3.4 °C 3:43 TU'
$ws.Range('O4').Value = '4.7 °C'
$ws.Range('E5').Value = '2026-02-09 04:18:30'
$ws.Range('H5').Value = '90%'
$ws.Range('O5').Value = '-5.6 °C'
$ws.Range('E6').Value = '2026-02-09 04:18:32'
$ws.Range('N6').Value = '5.3 °C 3:43 TU'
$ws.Range('O6').Value = '6.8 °C'
$ws.Range('E7').Value = '2026-02-09 04:18:34'
$ws.Range('H7').Value = '73%'
$ws.Range('E8').Value = '2026-02-09 04:18:37'
$ws.Range('H8').Value = '81%'
$ws.Range('E9').Value = '2026-02-09 04:18:39'
$ws.Range('H9').Value = '83%'
$ws.Range('N9').Value = '4.5 °C 3:39 TU'
$ws.Range('O9').Value = '7.7 °C'
$ws.Range('E10').Value = '2026-02-09 04:18:42'
$ws.Range('O10').Value = '6.2 °C'
$ws.Range('E11').Value = '2026-02-09 04:18:44'
$ws.Range('H11').Value = '99%'
$ws.Range('M11').Value = '2.8 °C 3:55 TU'
$ws.Range('O11').Value = '2.0 °C'
$ws.Range('E12').Value = '2026-02-09 04:18:47'
$ws.Range('H12').Value = '88%'
$ws.Range('O12').Value = '8.0 °C'
$ws.Range('E13').Value = '2026-02-09 04:18:49'
$ws.Range('H13').Value = '96%'
$ws.Range('J13').Value = '1011.1 hPa'
$ws.Range('N13').Value = '-3.2 °C 3:59 TU'
$ws.Range('O13').Value = '-1.4 °C'
$ws.Range('E14').Value = '2026-02-09 04:18:51'
$ws.Range('L14').Value = '13.0 km/h - 306º 3:50 TU'
$ws.Range('E15').Value = '2026-02-09 04:18:54'
$ws.Range('H15').Value = '85%'
$ws.Range('N15').Value = '3.2 °C 3:55 TU'
$ws.Range('O15').Value = '6.2 °C'
$ws.Range('E16').Value = '2026-02-09 04:18:56'
$ws.Range('H16').Value = '62%'
$ws.Range('M16').Value = '-4.3 °C 3:34 TU'
$ws.Range('O16').Value = '-5.1 °C'
$ws.Range('E17').Value = '2026-02-09 04:18:59'
$ws.Range('H17').Value = '94%'
$ws.Range('N17').Value = '-0.7 °C 3:30 TU'
$ws.Range('E18').Value = '2026-02-09 04:19:01'
$ws.Range('N18').Value = '4.8 °C 3:59 TU'
$ws.Range('O18').Value = '6.8 °C'
$ws.Range('E19').Value = '2026-02-09 04:19:04'
$ws.Range('L19').Value = '8.6 km/h - 230º 3:38 TU'
$ws.Range('N19').Value = '3.0 °C 3:58 TU'
$ws.Range('O19').Value = '3.3 °C'
$ws.Range('E20').Value = '2026-02-09 04:19:06'
$ws.Range('N20').Value = '-7.7 °C 3:31 TU'
$ws.Range('O20').Value = '-6.4 °C'
$ws.Range('E21').Value = '2026-02-09 04:19:08'
$ws.Range('H21').Value = '93%'
$ws.Range('J21').Value = '1009.8 hPa'
$ws.Range('N21').Value = '-0.4 °C 3:30 TU'
$ws.Range('O21').Value = '0.8 °C'
$ws.Range('E22').Value = '2026-02-09 04:19:11'
$ws.Range('M22').Value = '-7.2 °C 3:55 TU'
$ws.Range('O22').Value = '-7.6 °C'
$ws.Range('E23').Value = '2026-02-09 04:19:13'
$ws.Range('O23').Value = '-5.9 °C'
$ws.Range('E24').Value = '2026-02-09 04:19:16'
$ws.Range('H24').Value = '85%'
$ws.Range('N24').Value = '2.8 °C 3:56 TU'
$ws.Range('O24').Value = '4.6 °C'
$ws.Range('E25').Value = '2026-02-09 04:19:18'
$ws.Range('H25').Value = '78%'
$ws.Range('O25').Value = '-4.5 °C'
$ws.Range('E26').Value = '2026-02-09 04:19:21'
$ws.Range('H26').Value = '92%'
$ws.Range('N26').Value = '-0.8 °C 3:49 TU'
$ws.Range('E27').Value = '2026-02-09 04:19:23'
$ws.Range('N27').Value = '-4.5 °C 3:36 TU'
$ws.Range('E28').Value = '2026-02-09 04:19:26'
$ws.Range('N28').Value = '2.0 °C 3:56 TU'
$ws.Range('O28').Value = '3.9 °C'
$ws.Range('E29').Value = '2026-02-09 04:19:28'
$ws.Range('H29').Value = '95%'
$ws.Range('N29').Value = '3.8 °C 3:59 TU'
$ws.Range('O29').Value = '6.0 °C'
$ws.Range('E30').Value = '2026-02-09 04:19:31'
$ws.Range('J30').Value = '1007.8 hPa'
$ws.Range('N30').Value = '5.4 °C 3:57 TU'
$ws.Range('O30').Value = '7.0 °C'
$ws.Range('E31').Value = '2026-02-09 04:19:33'
$ws.Range('H31').Value = '73%'
$ws.Range('M31').Value = '10.1 °C 3:40 TU'
$ws.Range('O31').Value = '8.8 °C'
$ws.Range('E32').Value = '2026-02-09 04:19:35'
$ws.Range('H32').Value = '79%'
$ws.Range('N32').Value = '2.6 °C 3:34 TU'
$ws.Range('O32').Value = '3.2 °C'
$ws.Range('E33').Value = '2026-02-09 04:19:38'
$ws.Range('J33').Value = '1009.7 hPa'
$ws.Range('O33').Value = '-0.3 °C'
$ws.Range('E34').Value = '2026-02-09 04:19:40'
$ws.Range('N34').Value = '-5.2 °C 3:54 TU'
$ws.Range('O34').Value = '-2.8 °C'
$ws.Range('E35').Value = '2026-02-09 04:19:43'
$ws.Range('H35').Value = '67%'
$ws.Range('J35').Value = '1010.1 hPa'
$ws.Range('M35').Value = '4.5 °C 3:51 TU'
$ws.Range('O35').Value = '3.6 °C'
$ws.Range('E36').Value = '2026-02-09 04:19:45'
$ws.Range('H36').Value = '79%'
$ws.Range('J36').Value = '1007.7 hPa'
$ws.Range('O36').Value = '8.9 °C'
$ws.Range('E37').Value = '2026-02-09 04:19:48'
$ws.Range('N37').Value = '2.5 °C 3:59 TU'
$ws.Range('O37').Value = '3.8 °C'
$ws.Range('E38').Value = '2026-02-09 04:19:50'
$ws.Range('N38').Value = '5.4 °C 3:59 TU'
$ws.Range('E39').Value = '2026-02-09 04:19:52'
$ws.Range('E40').Value = '2026-02-09 04:19:55'
$ws.Range('J40').Value = '1010.7 hPa'
$ws.Range('N40').Value = '-1.1 °C 3:52 TU'
$ws.Range('O40').Value = '-0.2 °C'
$ws.Range('E41').Value = '2026-02-09 04:19:57'
$ws.Range('K41').Value = '-0.1 MJ/m2'
$ws.Range('E42').Value = '2026-02-09 04:19:59'
$ws.Range('H42').Value = '96%'
$ws.Range('N42').Value = '4.9 °C 3:51 TU'
$ws.Range('O42').Value = '6.8 °C'
$ws.Range('E43').Value = '2026-02-09 04:20:02'
$ws.Range('N43').Value = '5.9 °C 3:58 TU'
$ws.Range('O43').Value = '6.5 °C'
$ws.Range('E44').Value = '2026-02-09 04:20:04'
$ws.Range('E45').Value = '2026-02-09 04:20:07'
$ws.Range('N45').Value = '-0.9 °C 3:54 TU'
$ws.Range('O45').Value = '0.1 °C'
$ws.Range('E46').Value = '2026-02-09 04:20:09'
$ws.Range('H46').Value = '82%'
$ws.Range('J46').Value = '1009.6 hPa'
